# Auto-generated COM-interop edit script.
# Rerun-and-resummarise carown_LR/Magdeburg model outputs:
#  - rename/reorder the 7 existing 'summX' sheets
#  - collapse HHType_simp Single_Female_Parent/Single_Male_Parent
#    into a single Single_Parent category (one fewer row per sheet)
#  - write the freshly rerun coefficient/p-value table into every sheet
#  - append a brand-new 8th sheet ('summ9') with its own rerun results

$wb = $excel.ActiveWorkbook

# --- Step 1: stage renames through unique temp names (avoids collisions
#     from the summ17<->summ3 swap and the summ13->summ0->summ19 chain) ---
for ($i = 1; $i -le 7; $i++) {
    $wb.Worksheets.Item($i).Name = "zzz_tmp_${i}"
}

$finalNames = @("summ1", "summ19", "summ3", "summ17", "summ0", "summ2", "summ6")
for ($i = 1; $i -le 7; $i++) {
    $wb.Worksheets.Item($i).Name = $finalNames[$i - 1]
}

# --- Step 2: for each of the 7 original sheets, drop the row for the
#     category merged away (old row 8), relabel rows 6-7, and write
#     the freshly rerun coefficient/p-value numbers ---
# sheet 1: summ1
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -17.49237464737675
$data[0,1] = 0.9961025167630034
$data[1,0] = 17.82488037657322
$data[1,1] = 0.9960284304267272
$data[2,0] = 18.75161281503001
$data[2,1] = 0.9958219466309539
$data[3,0] = 15.72333079738496
$data[3,1] = 0.9964966744217983
$data[4,0] = 16.18896351532735
$data[4,1] = 0.9963929271881916
$data[5,0] = 17.91170697534538
$data[5,1] = 0.9960090847533692
$data[6,0] = -0.5207994106936459
$data[6,1] = 0.00498328290113684
$data[7,0] = 0.0003364866071974766
$data[7,1] = 0.0004176184226060882
$data[8,0] = -0.01651266922586832
$data[8,1] = 0.03412548558619004
$data[9,0] = -0.1250329067919562
$data[9,1] = 0.4890428993817817
$data[10,0] = 0.8012398103078283
$data[10,1] = 0.008749398308803114
$data[11,0] = 0.8613416966262368
$data[11,1] = 0.02555461823748297
$data[12,0] = -0.0003471599047933904
$data[12,1] = 0.08703685259983698
$data[13,0] = [double]"1.077411673830398e-07"
$data[13,1] = 0.4723714995556555
$data[14,0] = -0.1969759128276478
$data[14,1] = 0.4376146253072449
$data[15,0] = -0.004016777530475427
$data[15,1] = 0.991846255095909
$data[16,0] = -1.128502817111067
$data[16,1] = 0.8779623158644725
$data[17,0] = 0.06244296499710112
$data[17,1] = 0.07652225612476957
$data[18,0] = 0.009911591928487402
$data[18,1] = 0.5050036617052864
$data[19,0] = 4.376707103789719
$data[19,1] = 0.5771600715301966
$data[20,0] = 0.248777273467529
$data[20,1] = 0.9554786768302708
$data[21,0] = -1.940663626966046
$data[21,1] = 0.6749839810157026
$ws.Range("B2:C23").Value = $data

# sheet 2: summ19
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -15.51835052441645
$data[0,1] = 0.9973051390865679
$data[1,0] = 18.04359064109249
$data[1,1] = 0.996866616498112
$data[2,0] = 18.78637308398921
$data[2,1] = 0.9967376283574717
$data[3,0] = 16.04700471546651
$data[3,1] = 0.9972133347614407
$data[4,0] = 16.48891460037986
$data[4,1] = 0.9971365946015394
$data[5,0] = 18.06683689770867
$data[5,1] = 0.9968625796801355
$data[6,0] = -0.468205589676431
$data[6,1] = 0.0178728784119327
$data[7,0] = 0.0003894123892714915
$data[7,1] = [double]"9.745780732511054e-05"
$data[8,0] = -0.01805650785440854
$data[8,1] = 0.01939256642106701
$data[9,0] = -0.2109291744111666
$data[9,1] = 0.2503206254053305
$data[10,0] = 0.8406454268537845
$data[10,1] = 0.005574619516843981
$data[11,0] = 1.004309515064908
$data[11,1] = 0.008281961973347714
$data[12,0] = -0.0002776985193687094
$data[12,1] = 0.172759311844181
$data[13,0] = [double]"8.976485474909453e-09"
$data[13,1] = 0.9513657530366033
$data[14,0] = -0.2600617684517795
$data[14,1] = 0.3142333514611058
$data[15,0] = -0.01684273738773154
$data[15,1] = 0.9653235242407243
$data[16,0] = -0.8828123897046822
$data[16,1] = 0.9062887983381822
$data[17,0] = 0.02872737335638999
$data[17,1] = 0.4087548985068911
$data[18,0] = 0.002784463636143045
$data[18,1] = 0.8499099918619644
$data[19,0] = 5.169733331946808
$data[19,1] = 0.5031445672934152
$data[20,0] = 0.6271615311676961
$data[20,1] = 0.885026219175014
$data[21,0] = -2.237331159394627
$data[21,1] = 0.6221470399563167
$ws.Range("B2:C23").Value = $data

# sheet 3: summ3
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -17.66262888359355
$data[0,1] = 0.9969351542496225
$data[1,0] = 18.16330215918316
$data[1,1] = 0.9968482765308857
$data[2,0] = 19.11607489304614
$data[2,1] = 0.9966829509221211
$data[3,0] = 16.03027282540733
$data[3,1] = 0.9972184013869119
$data[4,0] = 16.54248780103299
$data[4,1] = 0.9971295214010087
$data[5,0] = 18.21138716824615
$data[5,1] = 0.9968399328111046
$data[6,0] = -0.5523065839287293
$data[6,1] = 0.004195668857314344
$data[7,0] = 0.0003451453387503962
$data[7,1] = 0.0005377439792426792
$data[8,0] = -0.01581440300223877
$data[8,1] = 0.0427199682492223
$data[9,0] = -0.1261223659344266
$data[9,1] = 0.4954695401212319
$data[10,0] = 0.969952152601713
$data[10,1] = 0.001289802512859197
$data[11,0] = 0.8674395482778209
$data[11,1] = 0.02205284230066775
$data[12,0] = -0.000306688402385918
$data[12,1] = 0.1434637513406948
$data[13,0] = [double]"9.524210735289296e-08"
$data[13,1] = 0.522458665034937
$data[14,0] = -0.1912938955444979
$data[14,1] = 0.4495752458681964
$data[15,0] = 0.1772335521138716
$data[15,1] = 0.6515836420131984
$data[16,0] = 5.179419597673303
$data[16,1] = 0.4898249413263994
$data[17,0] = 0.04355429880776716
$data[17,1] = 0.2212375637270571
$data[18,0] = 0.001005823537447946
$data[18,1] = 0.9456651196773335
$data[19,0] = 2.995899891441537
$data[19,1] = 0.7043957324243426
$data[20,0] = -0.6724058288817687
$data[20,1] = 0.8795333965449884
$data[21,0] = -0.1089336435545556
$data[21,1] = 0.9811711272992505
$ws.Range("B2:C23").Value = $data

# sheet 4: summ17
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -16.56885397033806
$data[0,1] = 0.997124469589928
$data[1,0] = 18.31082191638636
$data[1,1] = 0.9968221513823005
$data[2,0] = 19.07559480768115
$data[2,1] = 0.9966894255827902
$data[3,0] = 16.2419526093216
$data[3,1] = 0.9971812026816409
$data[4,0] = 16.58159374025687
$data[4,1] = 0.997122258081402
$data[5,0] = 18.01308621877696
$data[5,1] = 0.996873823251387
$data[6,0] = -0.5194992929072283
$data[6,1] = 0.007220137052808938
$data[7,0] = 0.0003678196139477162
$data[7,1] = 0.0002588235306731543
$data[8,0] = -0.01908448307647327
$data[8,1] = 0.01587845122294
$data[9,0] = -0.1562584816709722
$data[9,1] = 0.3952008064023799
$data[10,0] = 1.05959188521129
$data[10,1] = 0.0005075907689372026
$data[11,0] = 1.024964113943987
$data[11,1] = 0.007669552784451345
$data[12,0] = -0.0002544772622453354
$data[12,1] = 0.2220128790443934
$data[13,0] = [double]"5.751983043929519e-09"
$data[13,1] = 0.9684639154171732
$data[14,0] = -0.2040025531958419
$data[14,1] = 0.4245590677220712
$data[15,0] = -0.03122281929012949
$data[15,1] = 0.9362951202475908
$data[16,0] = 1.412895934569186
$data[16,1] = 0.8487481762998037
$data[17,0] = 0.02933081748027915
$data[17,1] = 0.3991829097493442
$data[18,0] = 0.007187505362661901
$data[18,1] = 0.6283427982004021
$data[19,0] = 2.481142561044861
$data[19,1] = 0.7541176005764498
$data[20,0] = -1.214977387602054
$data[20,1] = 0.7857728897372366
$data[21,0] = -0.5509328903235658
$data[21,1] = 0.9056999435722206
$ws.Range("B2:C23").Value = $data

# sheet 5: summ0
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -15.23841738766735
$data[0,1] = 0.9973545705104938
$data[1,0] = 18.44951073211951
$data[1,1] = 0.9967971183432249
$data[2,0] = 19.29968676551633
$data[2,1] = 0.9966495265030932
$data[3,0] = 16.11523212929111
$data[3,1] = 0.9972023532248148
$data[4,0] = 16.56050768469908
$data[4,1] = 0.9971250525440014
$data[5,0] = 18.33862057633862
$data[5,1] = 0.9968163690712626
$data[6,0] = -0.5968894781492099
$data[6,1] = 0.004759126313869502
$data[7,0] = 0.0003160398028654458
$data[7,1] = 0.001155960896411034
$data[8,0] = -0.01902526834178749
$data[8,1] = 0.01802230097916316
$data[9,0] = -0.1475052376685558
$data[9,1] = 0.4235519704736941
$data[10,0] = 0.9856244980339749
$data[10,1] = 0.00143564753830049
$data[11,0] = 0.934280394887481
$data[11,1] = 0.01868541171186711
$data[12,0] = -0.0002869784376477235
$data[12,1] = 0.1654638791852394
$data[13,0] = [double]"-5.494421250974823e-10"
$data[13,1] = 0.9970699240266255
$data[14,0] = -0.2480477520844884
$data[14,1] = 0.3331099615293132
$data[15,0] = -0.1097190359328545
$data[15,1] = 0.7789135025276854
$data[16,0] = -2.614576990144448
$data[16,1] = 0.7238944035732953
$data[17,0] = 0.03379600616712061
$data[17,1] = 0.3405272669248472
$data[18,0] = 0.005951627674896283
$data[18,1] = 0.6852699541401448
$data[19,0] = 6.856701576715877
$data[19,1] = 0.3794958152534176
$data[20,0] = 1.254022719634912
$data[20,1] = 0.7751728626308549
$data[21,0] = -3.500572639952563
$data[21,1] = 0.4447810050148625
$ws.Range("B2:C23").Value = $data

# sheet 6: summ2
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -14.14208010003364
$data[0,1] = 0.997543288796821
$data[1,0] = 18.09306827669193
$data[1,1] = 0.9968569398411148
$data[2,0] = 19.08455710454611
$data[2,1] = 0.9966847030966314
$data[3,0] = 16.04074839213003
$data[3,1] = 0.9972134596365573
$data[4,0] = 16.56007467859212
$data[4,1] = 0.9971232445508281
$data[5,0] = 18.17528420304119
$data[5,1] = 0.9968426576924835
$data[6,0] = -0.4375288355247033
$data[6,1] = 0.03963942218182009
$data[7,0] = 0.0004202635609304073
$data[7,1] = [double]"6.317599816316764e-05"
$data[8,0] = -0.02219871970787233
$data[8,1] = 0.005335538347488063
$data[9,0] = -0.109707356426302
$data[9,1] = 0.5531344949373997
$data[10,0] = 0.9800283174839328
$data[10,1] = 0.001600757131137975
$data[11,0] = 1.121269294773751
$data[11,1] = 0.004277326404708791
$data[12,0] = -0.0003490855459232983
$data[12,1] = 0.0852162734672947
$data[13,0] = [double]"-6.947527084937284e-08"
$data[13,1] = 0.6594167298783371
$data[14,0] = -0.3742249157815378
$data[14,1] = 0.1756763699219961
$data[15,0] = -0.1821811226897376
$data[15,1] = 0.6495434636822974
$data[16,0] = -4.187905577590018
$data[16,1] = 0.5936536476913904
$data[17,0] = 0.01396683146095513
$data[17,1] = 0.7004764552574416
$data[18,0] = 0.008542263358080244
$data[18,1] = 0.5677085727549892
$data[19,0] = 7.945073683841726
$data[19,1] = 0.3119324910944342
$data[20,0] = 1.552258041053323
$data[20,1] = 0.7224921814256567
$data[21,0] = -3.882209783662996
$data[21,1] = 0.4052062888793474
$ws.Range("B2:C23").Value = $data

# sheet 7: summ6
$ws = $wb.Worksheets.Item(7)
$ws.Rows.Item(8).Delete()
$ws.Range("A6").Value = "HHType_simp[T.Single_Male]"
$ws.Range("A7").Value = "HHType_simp[T.Single_Parent]"
$data = New-Object 'object[,]' 22,2
$data[0,0] = -17.58588759110744
$data[0,1] = 0.9969492888798097
$data[1,0] = 18.41800440936983
$data[1,1] = 0.9968049375923728
$data[2,0] = 19.23140299762627
$data[2,1] = 0.9966638341336855
$data[3,0] = 16.08090496587117
$data[3,1] = 0.9972103640328194
$data[4,0] = 16.49880516343119
$data[4,1] = 0.9971378690769932
$data[5,0] = 18.21442276488261
$data[5,1] = 0.9968402537392812
$data[6,0] = -0.6036563099175701
$data[6,1] = 0.001798590293084692
$data[7,0] = 0.0003400952498123117
$data[7,1] = 0.0007212620558648027
$data[8,0] = -0.01729745849372905
$data[8,1] = 0.03052560771470325
$data[9,0] = -0.235958744922774
$data[9,1] = 0.2056285686073331
$data[10,0] = 1.032519282567061
$data[10,1] = 0.0008064232348138526
$data[11,0] = 1.011190734960453
$data[11,1] = 0.00951196998189579
$data[12,0] = -0.0003582882879062629
$data[12,1] = 0.08864251268984924
$data[13,0] = [double]"7.32125186307272e-08"
$data[13,1] = 0.6316621168942582
$data[14,0] = -0.2105070299865066
$data[14,1] = 0.4288687457931395
$data[15,0] = -0.01217569607135187
$data[15,1] = 0.976199014599839
$data[16,0] = -0.3175546293943738
$data[16,1] = 0.9662102437846956
$data[17,0] = 0.06329693092958721
$data[17,1] = 0.0771492987631576
$data[18,0] = 0.009121470058173485
$data[18,1] = 0.5612590088477077
$data[19,0] = 4.860287412473141
$data[19,1] = 0.5565714838684497
$data[20,0] = 1.02229682479006
$data[20,1] = 0.827216130850955
$data[21,0] = -2.247997474408195
$data[21,1] = 0.6439163768913699
$ws.Range("B2:C23").Value = $data

# --- Step 3: append the brand-new 8th sheet ('summ9') ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "summ9"
$newSheet.Range("A1").Value = "param"
$newSheet.Range("B1").Value = "coefficient"
$newSheet.Range("C1").Value = "p"
$labels8 = @("Intercept", "HHType_simp[T.MultiAdult]", "HHType_simp[T.MultiAdult_Kids]", "HHType_simp[T.Single_Female]", "HHType_simp[T.Single_Male]", "HHType_simp[T.Single_Parent]", "HHSize", "IncomeDetailed_Numeric", "maxAgeHH", "UniversityEducation", "InEmployment", "AllRetired", "UrbPopDensity", "UrbBuildDensity", "DistSubcenter", "DistCenter", "bike_lane_share", "IntersecDensity", "StreetLength", "LU_UrbFab", "LU_Comm", "LU_Urban")
for ($i = 0; $i -lt 22; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $labels8[$i]
}
$data8 = New-Object 'object[,]' 22,2
$data8[0,0] = -17.58296816377954
$data8[0,1] = 0.9969497318726802
$data8[1,0] = 18.43637994643299
$data8[1,1] = 0.9968016834113999
$data8[2,0] = 19.20110773326211
$data8[2,1] = 0.9966690203116257
$data8[3,0] = 16.37389051601846
$data8[3,1] = 0.9971594794972052
$data8[4,0] = 16.80480809386572
$data8[4,1] = 0.9970847248231871
$data8[5,0] = 18.15905754382795
$data8[5,1] = 0.9968497926670734
$data8[6,0] = -0.443062857663666
$data8[6,1] = 0.03456208839695292
$data8[7,0] = 0.0003264584695579784
$data8[7,1] = 0.0009207445334723695
$data8[8,0] = -0.02028465050151514
$data8[8,1] = 0.01065957223884458
$data8[9,0] = -0.1996011928032508
$data8[9,1] = 0.2775254385765263
$data8[10,0] = 0.9458985169646225
$data8[10,1] = 0.002570793695654447
$data8[11,0] = 1.06474099537618
$data8[11,1] = 0.006036873884041694
$data8[12,0] = -0.0003134304064567499
$data8[12,1] = 0.1272163741520661
$data8[13,0] = [double]"7.910027144596754e-08"
$data8[13,1] = 0.5886560741441116
$data8[14,0] = -0.1066879541245566
$data8[14,1] = 0.6687730254546163
$data8[15,0] = 0.06694613109356823
$data8[15,1] = 0.8619119621810907
$data8[16,0] = 1.85331528674608
$data8[16,1] = 0.8018997073251838
$data8[17,0] = 0.04538846909814803
$data8[17,1] = 0.1968560632228703
$data8[18,0] = 0.003348777112423134
$data8[18,1] = 0.8174721214553478
$data8[19,0] = 1.366764412747428
$data8[19,1] = 0.8590359141621099
$data8[20,0] = -1.52960452930282
$data8[20,1] = 0.7224648684868349
$data8[21,0] = 0.3339323371790238
$data8[21,1] = 0.9405066697114386
$newSheet.Range("B2:C23").Value = $data8

# Match page margins (inches) used by the other sheets: 0.75/0.75/1/1/0.5/0.5
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

Write-Output "done"
